# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.984.01"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "1.818.90"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.20"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4304"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3693"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07256"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "2.125.08"
$ws.Range("E10").Value = "  +21.29%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8702"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.35"
$ws.Range("E12").Value = "  +5.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.415"
$ws.Range("E13").Value = "  +3.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.631"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06942"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.26"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.012"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008908"
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.26"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").Value = "27.027.50"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.211"
$ws.Range("E22").Value = "  +3.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.03"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").Value = "2.371.94"
$ws.Range("E24").Value = "  +19.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.07"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.888"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.42"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.251"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.932"
$ws.Range("E29").Value = "  +13.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.95"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08971"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.189"
$ws.Range("E32").Value = "  +6.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7519"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.429"
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.806"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.005"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.134"
$ws.Range("E37").Value = "  +5.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05246"
$ws.Range("E38").Value = "  +2.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01930"
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5128"
$ws.Range("E40").Value = "  +4.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1657"
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.744"
$ws.Range("E42").Value = "  +8.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.486"
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.346"
$ws.Range("E44").Value = "  +3.33%  "
$ws.Range("B45").Value = "PaxosStandard"
$ws.Range("C45").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.007"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.28"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.39"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4581"
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.653"
$ws.Range("E50").Value = "  +4.49%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06220"
$ws.Range("E51").Value = "  +0.50%  "
